# Apply the "ages" workbook edits:
#  - clear Tomas's age (B6) on the "ages" sheet
#  - add a new row for Bob (age 45) at row 9 on the "ages" sheet
#  - make "ages" the active sheet/tab, with F12 selected
#  - leave "Sheet1" selection at C16 (already the case) but no longer the active tab

$wb = $excel.ActiveWorkbook

$wsAges = $wb.Worksheets.Item("ages")
$wsSheet1 = $wb.Worksheets.Item("Sheet1")

# Clear Tomas's age value (B6)
$wsAges.Range("B6").ClearContents()

# Add the new Bob row (row 8 is intentionally left blank, data resumes at row 9)
$wsAges.Range("A9").Value = "Bob"
$wsAges.Range("B9").Value = 45

# Keep Sheet1's existing selection
$wsSheet1.Range("C16").Select()

# Activate the "ages" sheet and select F12, making it the visible/active tab
$wsAges.Activate()
$wsAges.Range("F12").Select()
